$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.577.53'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.921.70'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.49'
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4837'
$ws.Range("E7").Value = '  +2.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2904'
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06804'
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '112.12'
$ws.Range("E10").Value = '  +6.17%  '

$ws.Range("E11").Value = '  +5.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.916.12'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.484'
$ws.Range("E13").Value = '  +2.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07579'
$ws.Range("E14").Value = '  -1.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6748'
$ws.Range("E15").Value = '  +0.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '294.40'
$ws.Range("E16").Value = '  +1.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.570.96'
$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007687'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.05'
$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9996'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.529'
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.158.57'
$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.468'
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.505'
$ws.Range("E25").Value = '  -0.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.14'
$ws.Range("E26").Value = '  -0.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.36'
$ws.Range("E27").Value = '  -2.40%  '

$ws.Range("E28").Value = '  -1.01%  '

$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("E30").Value = '  +2.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.143'
$ws.Range("E31").Value = '  -0.90%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.069'
$ws.Range("E32").Value = '  +0.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04963'
$ws.Range("E33").Value = '  -1.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7363'
$ws.Range("E34").Value = '  +0.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.141'
$ws.Range("E35").Value = '  -0.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02034'
$ws.Range("E36").Value = '  -1.69%  '

$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.689'
$ws.Range("E38").Value = '  +0.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.026'
$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '109.72'
$ws.Range("E40").Value = '  -1.58%  '

$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8716'
$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.859'
$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.38'
$ws.Range("E45").Value = '  +2.37%  '

$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.88'
$ws.Range("E47").Value = '  +1.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.234'
$ws.Range("E48").Value = '  -1.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1231'
$ws.Range("E49").Value = '  -1.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.91'
$ws.Range("E50").Value = '  -0.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2512'
$ws.Range("E51").Value = '  -0.04%  '
